$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Addl Shortcuts")

$cmd = [char]0x2318

# Rows 16-18: capitalize 'Control' -> 'CONTROL' in the keyboard-shortcut text
$controlText = "$cmd click 'CONTROL' Button [Mac]\nCTRL click 'CONTROL' Button [Win]"
$ws.Range("A16").Value = $controlText
$ws.Range("A17").Value = $controlText
$ws.Range("A18").Value = $controlText

# Insert a new row 19 documenting the Graph Resolution % shortcut
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "$cmd+PLUS, $cmd+MINUS [Mac]\nCTRL+SHIFT+PLUS, CTRL+MINUS [Win]"
$ws.Range("B19").Value = "Graph"
$ws.Range("C19").Value = "Increase or Decrease Graph Resolution %"
$ws.Range("D19").Value = "Same as Config>> Curves> UI tab>> Graph % +/-"
